$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("BQ2").Value = 86
$ws.Range("BU2").Value = 207.1
$ws.Range("BW2").Value = 321
$ws.Range("CA2").Value = 46
$ws.Range("CB2").Value = 40
$ws.Range("CE2").Value = 2.08
$ws.Range("CH2").Value = 7.82
$ws.Range("CJ2").Value = 18.83
$ws.Range("CL2").Value = 29.18
$ws.Range("CP2").Value = 4.18
$ws.Range("CQ2").Value = 3.64

# Row 3
$ws.Range("BQ3").Value = 71
$ws.Range("BU3").Value = 253.09
$ws.Range("CH3").Value = 6.45
$ws.Range("CJ3").Value = 23.01

# Row 4
$ws.Range("BW4").Value = 225
$ws.Range("BY4").Value = 2
$ws.Range("BZ4").Value = 118
$ws.Range("CA4").Value = 0
$ws.Range("CB4").Value = 16
$ws.Range("CC4").Value = 35
$ws.Range("CD4").Value = 54
$ws.Range("CE4").Value = 1.51
$ws.Range("CL4").Value = 20.45
$ws.Range("CN4").Value = 0.18
$ws.Range("CO4").Value = 10.73
$ws.Range("CP4").Value = 0
$ws.Range("CQ4").Value = 1.45
$ws.Range("CR4").Value = 3.18
$ws.Range("CS4").Value = 4.91

# Row 6
$ws.Range("BQ6").Value = 80
$ws.Range("BR6").Value = 196
$ws.Range("BU6").Value = 276.9
$ws.Range("BV6").Value = 39.31
$ws.Range("BW6").Value = 291
$ws.Range("BY6").Value = 19
$ws.Range("CE6").Value = 2.21
$ws.Range("CH6").Value = 6.67
$ws.Range("CI6").Value = 16.33
$ws.Range("CJ6").Value = 23.07
$ws.Range("CK6").Value = 3.28
$ws.Range("CL6").Value = 24.25
$ws.Range("CN6").Value = 1.58

# Row 7
$ws.Range("BQ7").Value = 66
$ws.Range("BR7").Value = 113
$ws.Range("BU7").Value = 186.78
$ws.Range("BV7").Value = 34.42
$ws.Range("BW7").Value = 186
$ws.Range("CB7").Value = 12
$ws.Range("CE7").Value = 2.5
$ws.Range("CH7").Value = 6
$ws.Range("CI7").Value = 10.27
$ws.Range("CJ7").Value = 16.98
$ws.Range("CK7").Value = 3.13
$ws.Range("CL7").Value = 16.91
$ws.Range("CQ7").Value = 1.09

# Row 8
$ws.Range("BR8").Value = 273
$ws.Range("BV8").Value = 44.36
$ws.Range("BW8").Value = 348
$ws.Range("CA8").Value = 4
$ws.Range("CB8").Value = 12
$ws.Range("CE8").Value = 1.79
$ws.Range("CI8").Value = 22.75
$ws.Range("CK8").Value = 3.7
$ws.Range("CL8").Value = 29
$ws.Range("CP8").Value = 0.33
$ws.Range("CQ8").Value = 1
